# Fruta / hortaliza, semanal
# Insert a new weekly record at row 11 (pushing the existing rows 11-26 down to 12-27)
# for "Vega Modelo de Temuco - Maracuyá".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 11, shifting existing data down
# (mirrors the rows below/above so formatting such as the date style on
# column D is carried over to the new row, matching the target workbook).
$ws.Rows.Item(11).Insert()

$ws.Cells.Item(11, 1).Value  = 10
$ws.Cells.Item(11, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(11, 3).Value  = "La Araucanía"
$ws.Cells.Item(11, 4).Value  = 44645
$ws.Cells.Item(11, 5).Value  = 9
$ws.Cells.Item(11, 6).Value  = "Fruta"
$ws.Cells.Item(11, 7).Value  = 100108
$ws.Cells.Item(11, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(11, 9).Value  = 100108003
$ws.Cells.Item(11, 10).Value = "Maracuyá"
$ws.Cells.Item(11, 11).Value = "Sin especificar"
$ws.Cells.Item(11, 12).Value = "Primera"
$ws.Cells.Item(11, 13).Value = 5
$ws.Cells.Item(11, 14).Value = 30000
$ws.Cells.Item(11, 15).Value = 30000
$ws.Cells.Item(11, 16).Value = 30000
$ws.Cells.Item(11, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(11, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(11, 19).Value = 1667
$ws.Cells.Item(11, 20).Value = 18
